$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend formatting down to the two new rows (17-18) before writing data ---
# Copy the formatting (styles + row height) of the last existing data row onto the two new rows.
$ws.Range("B16:F16").Copy() | Out-Null
$ws.Range("B17:F18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Rows("17:18").RowHeight = $ws.Rows(16).RowHeight

# --- Clear every existing hyperlink up front; row positions are about to shift, so remapping
#     ref->rId cleanly after the fact is simpler than trying to preserve the old ones. ---
$ws.Hyperlinks.Delete()

# --- Rewrite every data row (3-18) with the final values from the edited dataset ---
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 43522
$ws.Range("D3").Value = 'https://www.quora.com/Can-you-collect-unemployment-if-you-cant-get-a-job-post-graduation'
$ws.Range("E3").Value = 'Can You Collect Unemployment if You Can''t Get a Job Post-Graduation?'
$ws.Range("F3").Value = 'The short answer is “No”, but there is a longer answer that may lead you to have some hope…'

$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 43004
$ws.Range("D4").Value = 'https://bizfluent.com/info-10060450-can-college-students-receive-unemployment-benefits.html'
$ws.Range("E4").Value = 'Can College Students Receive Unemployment Benefits?'
$ws.Range("F4").Value = 'The traditional image of the unemployed doesn’t usually include college students. However, many college students …'

$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 43420
$ws.Range("D5").Value = 'https://www.thebalancecareers.com/is-private-unemployment-insurance-worth-it-4161288'
$ws.Range("E5").Value = 'What Private Unemployment Insurance Is Available? '
$ws.Range("F5").Value = 'What is private unemployment insurance? Is it worth it? Learn where you can get job loss insurance coverage …'

$ws.Range("B6").Value = 4
$ws.Range("C6").Value = 43524
$ws.Range("D6").Value = 'https://www.quora.com/Is-it-normal-for-a-college-student-to-be-scared-about-the-future'
$ws.Range("E6").Value = 'Is it Normal for a College Student to be Scared About the Future?'
$ws.Range("F6").Value = 'Honestly, it is normal for everyone to be scared about the future; that includes college students…'

$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 42991
$ws.Range("D7").Value = 'https://www.usnews.com/education/best-colleges/articles/2017-09-13/find-out-where-college-graduates-go-after-graduation'
$ws.Range("E7").Value = 'Find Out Where College Graduates Go After Graduation'
$ws.Range("F7").Value = 'Based on U.S. News data, college graduates choose graduate education programs over law or medical school.'

$ws.Range("B8").Value = 6
$ws.Range("C8").Value = 43525
$ws.Range("D8").Value = 'https://www.quora.com/Why-are-so-many-young-graduates-unemployed-or-underemployed'
$ws.Range("E8").Value = 'Why Are So Many Young Graduates Unemployed or Underemployed?'
$ws.Range("F8").Value = 'Although this feels very true, the available data in the U.S. doesn’t really look as bad as one would think…'

$ws.Range("B9").Value = 7
$ws.Range("C9").Value = 42677
$ws.Range("D9").Value = 'https://www.forbes.com/sites/realspin/2016/11/03/solving-the-college-affordability-problem-with-student-loan-insurance/'
$ws.Range("E9").Value = 'Solving the College Affordability Problem with Student Loan Insurance'
$ws.Range("F9").Value = 'Taxpayer dollars are scarce. There are alternative ways to keep college affordable, but doing so requires rethinking …'

$ws.Range("B10").Value = 8
$ws.Range("C10").Value = 42590
$ws.Range("D10").Value = 'https://studentloanhero.com/featured/unemployed-college-graduates-jobless/'
$ws.Range("E10").Value = 'What to Do If You Can’t Find a Job After Graduation'
$ws.Range("F10").Value = 'Calling all unemployed college graduates — here''s how to spend less and earn more if you''re struggling to find a job …'

$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 43521
$ws.Range("D11").Value = 'https://www.quora.com/Will-I-be-able-to-get-a-job-after-college'
$ws.Range("E11").Value = 'Will I Be Able to Get a Job After College?'
$ws.Range("F11").Value = 'Economically speaking, it’s normal to wonder if the “opportunity cost” of getting an education is going to be worth the payoff…'

$ws.Range("B12").Value = 10
$ws.Range("C12").Value = 42538
$ws.Range("D12").Value = 'https://research.stlouisfed.org/publications/review/2016/06/17/student-loans-under-the-risk-of-youth-unemployment/'
$ws.Range("E12").Value = 'Student Loans Under the Risk of Youth Unemployment '
$ws.Range("F12").Value = 'While most college graduates eventually find jobs that match their qualifications, the possibility of long spells of unemployment …'

$ws.Range("B13").Value = 11
$ws.Range("C13").Value = 43520
$ws.Range("D13").Value = 'https://www.quora.com/What-is-the-unemployment-rate-for-recent-college-grads'
$ws.Range("E13").Value = 'What is the Unemployment Rate for Recent College Grads?'
$ws.Range("F13").Value = 'In the U.S., there’s actually a whole division of the Department of Education devoted to tracking information like this…'

$ws.Range("B14").Value = 12
$ws.Range("C14").Value = 42518
$ws.Range("D14").Value = 'https://www.nytimes.com/2016/05/28/your-money/finally-private-unemployment-insurance-but-will-anyone-buy-it.html'
$ws.Range("E14").Value = 'Finally, Private Unemployment Insurance. But Will Anyone Buy It?'
$ws.Range("F14").Value = 'IncomeAssure offers coverage that provides half your pretax pay, minus state jobless benefits, for up to 24 weeks.'

$ws.Range("B15").Value = 13
$ws.Range("C15").Value = 40773
$ws.Range("D15").Value = 'https://www.sapling.com/11368489/can-collect-unemployment-after-graduation-graduate-school'
$ws.Range("E15").Value = 'Can I Collect Unemployment After Graduation From Graduate School?'
$ws.Range("F15").Value = 'Even if you applied yourself diligently throughout your final year in grad school seeking a job to kick off your career …'

$ws.Range("B16").Value = 14
$ws.Range("C16").Value = 43523
$ws.Range("D16").Value = 'https://www.quora.com/Who-pays-unemployment-insurance'
$ws.Range("E16").Value = 'Who Pays Unemployment Insurance?'
$ws.Range("F16").Value = 'In the U.S. it depends on who is offering the unemployment insurance. There are two main types of unemployment insurance…'

$ws.Range("B17").Value = 15
$ws.Range("C17").Value = 41334
$ws.Range("D17").Value = 'https://www.moneycrashers.com/cant-find-job-after-college/'
$ws.Range("E17").Value = 'What to Do If You Can’t Find a Job After College Graduation'
$ws.Range("F17").Value = 'Struggling to get a real job after college? See these crucial tips to better enjoy life and increase your chances of landing that first job.'

$ws.Range("B18").Value = 16
$ws.Range("C18").Value = 40033
$ws.Range("D18").Value = 'https://www.nytimes.com/2009/08/08/your-money/08money.html'
$ws.Range("E18").Value = 'Good Luck Getting Private Insurance for Unemployment'
$ws.Range("F18").Value = 'Paycheck replacement coverage is scarce, if it exists at all. …'

# --- Re-create the hyperlinks on column D for every row, pointing at the (now correct) URL ---
$ws.Hyperlinks.Add($ws.Range("D3"), 'https://www.quora.com/Can-you-collect-unemployment-if-you-cant-get-a-job-post-graduation') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), 'https://bizfluent.com/info-10060450-can-college-students-receive-unemployment-benefits.html') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), 'https://www.thebalancecareers.com/is-private-unemployment-insurance-worth-it-4161288') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D6"), 'https://www.quora.com/Is-it-normal-for-a-college-student-to-be-scared-about-the-future') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D7"), 'https://www.usnews.com/education/best-colleges/articles/2017-09-13/find-out-where-college-graduates-go-after-graduation') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D8"), 'https://www.quora.com/Why-are-so-many-young-graduates-unemployed-or-underemployed') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D9"), 'https://www.forbes.com/sites/realspin/2016/11/03/solving-the-college-affordability-problem-with-student-loan-insurance/') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D10"), 'https://studentloanhero.com/featured/unemployed-college-graduates-jobless/') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D11"), 'https://www.quora.com/Will-I-be-able-to-get-a-job-after-college') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D12"), 'https://research.stlouisfed.org/publications/review/2016/06/17/student-loans-under-the-risk-of-youth-unemployment/') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D13"), 'https://www.quora.com/What-is-the-unemployment-rate-for-recent-college-grads') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D14"), 'https://www.nytimes.com/2016/05/28/your-money/finally-private-unemployment-insurance-but-will-anyone-buy-it.html') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D15"), 'https://www.sapling.com/11368489/can-collect-unemployment-after-graduation-graduate-school') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D16"), 'https://www.quora.com/Who-pays-unemployment-insurance') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D17"), 'https://www.moneycrashers.com/cant-find-job-after-college/') | Out-Null
$ws.Hyperlinks.Add($ws.Range("D18"), 'https://www.nytimes.com/2009/08/08/your-money/08money.html') | Out-Null

# --- Refresh the AutoFilter over the new data extent and re-sort by Order ascending ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B18")) | Out-Null
$ws.Sort.SetRange($ws.Range("B2:F18"))
$ws.Sort.Header = 1
$ws.Sort.Apply()
$ws.Range("B2:F18").AutoFilter() | Out-Null

# --- Point the hidden _FilterDatabase defined name at the new filter range ---
$fd = $wb.Names.Item("List!_FilterDatabase")
$fd.RefersTo = "=List!`$B`$2:`$F`$18"

# --- Restore the active selection to B3, matching the post-edit workbook state ---
$ws.Range("B3").Select()
